# update MCA data source
# Applies the "update MCA data source" edit:
#  - data_source sheet: rows 6-10, Data Source column switches from OSM/geofabrik
#    to GEE / gee-community-catalog sources (with some links cleared / moved),
#    and row 10 is cleared out entirely.
#  - MCA_criteria sheet: a new "Exclusive Range" column (I) is added, and a
#    "Vector" data-type is recorded for the Existing Windfarm row (row 7).

$wb = $excel.ActiveWorkbook

$dataSource = $wb.Worksheets.Item("data_source")
$mcaCriteria = $wb.Worksheets.Item("MCA_criteria")

# --- data_source sheet -----------------------------------------------------

# Row 6: Water Area -> Bathymetry, sourced from GEE's Globathy dataset
$dataSource.Range("A6").Value = "Bathymetry"
$dataSource.Range("B6").Value = "GEE"
$dataSource.Range("C6").Value = "https://gee-community-catalog.org/projects/globathy/?h=globath"

# Row 7: Residencial Area keeps its label, but the link is dropped (GEE, no URL)
$dataSource.Range("B7").Value = "GEE"
$dataSource.Range("C7").ClearContents()

# Row 8: Roads, now linked to the GEE GRIP roads dataset
$dataSource.Range("B8").Value = "GEE"
$dataSource.Range("C8").Value = "https://gee-community-catalog.org/projects/grip/?h=roads"

# Row 9: Existing Windfarm, now linked (with an actual hyperlink) to the GEE
# energy farms dataset
$dataSource.Range("B9").Value = "GEE"
$dataSource.Range("C9").Value = "https://gee-community-catalog.org/projects/energy_farms/?h=wind+farm"
$dataSource.Hyperlinks.Add($dataSource.Range("C9"), "https://gee-community-catalog.org/projects/energy_farms/?h=wind+farm") | Out-Null

# Row 10: Substation row is removed entirely
$dataSource.Range("A10:C10").ClearContents()

# --- MCA_criteria sheet -----------------------------------------------------

# New column I: Exclusive Range header
$mcaCriteria.Range("I1").Value = "Exclusive Range"
$mcaCriteria.Columns.Item(9).ColumnWidth = 14.6

# Row 7 (Existing Windfarm): record its data type as Vector
$mcaCriteria.Range("C7").Value = "Vector"

# --- restore cursor / selection position -----------------------------------
# (matches the saved selection state: MCA_criteria cursor on A10, data_source
# cursor on B7, with data_source left as the active/visible tab)
$mcaCriteria.Range("A10").Select() | Out-Null
$dataSource.Activate()
$dataSource.Range("B7").Select() | Out-Null
